# Appends a new appointment row (row 3) to the worksheet, matching the
# author's edit: Full Name / Mobile / Reason / (new) Date-Time columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Deep Chakravorty"
# Mobile numbers must stay text (leading apostrophe forces text entry,
# same as the existing Mobile column values) rather than being coerced
# into a numeric value.
$ws.Range("B3").Value = "'9875480108"
$ws.Range("C3").Value = "I don't know."
$ws.Range("D3").Value = "2025-09-21 13:35:00"

# Reset to the default/Normal style so the new cells don't pick up an
# extra "quote prefix" number format and stay visually identical to the
# rest of the sheet.
$ws.Range("A3:D3").Style = "Normal"
